$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 600
$ws.Range("J12").Value = 599.5
$ws.Range("L12").Value = 599.5
$ws.Range("N12").Value = -939.5
$ws.Range("H33").Value = 380.77274
$ws.Range("I33").Value = 249.38889
$ws.Range("K33").Value = 249.38889
$ws.Range("M33").Value = -20.38889
$ws.Range("H51").Value = 3102
$ws.Range("J51").Value = 3102
$ws.Range("L51").Value = 3102
$ws.Range("N51").Value = -4070
$ws.Range("H76").Value = 5131.0557
$ws.Range("I76").Value = 4728.077
$ws.Range("J76").Value = 6178.8
$ws.Range("K76").Value = 4728.077
$ws.Range("L76").Value = 6178.8
$ws.Range("M76").Value = -4413.077
$ws.Range("N76").Value = -6808.8
$ws.Range("H79").Value = 5131.0557
$ws.Range("I79").Value = 4728.077
$ws.Range("J79").Value = 6178.8
$ws.Range("K79").Value = 4728.077
$ws.Range("L79").Value = 6178.8
$ws.Range("M79").Value = -3636.077
$ws.Range("N79").Value = -8362.799999999999
$ws.Range("H100").Value = 998
$ws.Range("I100").Value = 997.6667
$ws.Range("K100").Value = 997.6667
$ws.Range("M100").Value = -456.6667
$ws.Range("H111").Value = 617
$ws.Range("I111").Value = 763
$ws.Range("K111").Value = 2289
$ws.Range("M111").Value = 778
$ws.Range("H113").Value = 5063.5
$ws.Range("J113").Value = 6114.7
$ws.Range("L113").Value = 6114.7
$ws.Range("N113").Value = -12622.7
$ws.Range("H125").Value = 1578.375
$ws.Range("I125").Value = 2167
$ws.Range("K125").Value = 19503
$ws.Range("M125").Value = -17043

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1505
$ws.Range("I102").Value = 1106
$ws.Range("K102").Value = 1106
$ws.Range("M102").Value = 516
$ws.Range("H110").Value = 1429.28
$ws.Range("I110").Value = 1358.9048
$ws.Range("J110").Value = 1798.75
$ws.Range("K110").Value = 1358.9048
$ws.Range("L110").Value = 1798.75
$ws.Range("M110").Value = 686.0952
$ws.Range("N110").Value = -5888.75
$ws.Range("H132").Value = 1413.6562
$ws.Range("I132").Value = 1366.069
$ws.Range("K132").Value = 4098.207
$ws.Range("M132").Value = -1568.207

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2642.6667
$ws.Range("I99").Value = 2642.6667
$ws.Range("K99").Value = 2642.6667
$ws.Range("M99").Value = -1144.6667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 8921.6
$ws.Range("I62").Value = 9331.666999999999
$ws.Range("J62").Value = 8306.5
$ws.Range("K62").Value = 9331.666999999999
$ws.Range("L62").Value = 8306.5
$ws.Range("M62").Value = -8707.666999999999
$ws.Range("N62").Value = -9554.5
$ws.Range("H65").Value = 8921.6
$ws.Range("I65").Value = 9331.666999999999
$ws.Range("J65").Value = 8306.5
$ws.Range("K65").Value = 46658.335
$ws.Range("L65").Value = 41532.5
$ws.Range("M65").Value = -43538.335
$ws.Range("N65").Value = -47772.5
$ws.Range("H99").Value = 2820
$ws.Range("I99").Value = 2820
$ws.Range("K99").Value = 2820
$ws.Range("M99").Value = -1322
$ws.Range("H126").Value = 2820
$ws.Range("I126").Value = 2820
$ws.Range("K126").Value = 8460
$ws.Range("M126").Value = -5990

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 278.5
$ws.Range("I7").Value = 338
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 1014
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -902
$ws.Range("N7").Value = -524

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11795.235
$ws.Range("J80").Value = 17161.334
$ws.Range("L80").Value = 17161.334
$ws.Range("N80").Value = -19157.334
$ws.Range("H83").Value = 11795.235
$ws.Range("J83").Value = 17161.334
$ws.Range("L83").Value = 85806.67
$ws.Range("N83").Value = -95790.67
$ws.Range("H97").Value = 2258
$ws.Range("I97").Value = 2450
$ws.Range("J97").Value = 1874
$ws.Range("K97").Value = 2450
$ws.Range("L97").Value = 1874
$ws.Range("M97").Value = -1954
$ws.Range("N97").Value = -2866
$ws.Range("H113").Value = 2898.2727
$ws.Range("I113").Value = 2240.4285
$ws.Range("K113").Value = 2240.4285
$ws.Range("M113").Value = -70.42849999999999
$ws.Range("H126").Value = 3946.8
$ws.Range("I126").Value = 3915
$ws.Range("J126").Value = 3994.5
$ws.Range("K126").Value = 11745
$ws.Range("L126").Value = 11983.5
$ws.Range("M126").Value = -9275
$ws.Range("N126").Value = -16923.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3648.6
$ws.Range("I7").Value = 2248
$ws.Range("J7").Value = 3998.75
$ws.Range("K7").Value = 2248
$ws.Range("L7").Value = 3998.75
$ws.Range("M7").Value = -2136
$ws.Range("N7").Value = -4222.75
$ws.Range("H40").Value = 4573.25
$ws.Range("I40").Value = 4497
$ws.Range("J40").Value = 4649.5
$ws.Range("K40").Value = 4497
$ws.Range("L40").Value = 4649.5
$ws.Range("M40").Value = -4361
$ws.Range("N40").Value = -4921.5
$ws.Range("H55").Value = 384.9524
$ws.Range("I55").Value = 443.05884
$ws.Range("J55").Value = 138
$ws.Range("K55").Value = 443.05884
$ws.Range("L55").Value = 138
$ws.Range("M55").Value = -270.05884
$ws.Range("N55").Value = -484
$ws.Range("H68").Value = 3167.6155
$ws.Range("J68").Value = 3073.25
$ws.Range("L68").Value = 3073.25
$ws.Range("N68").Value = -4571.25
$ws.Range("H71").Value = 3167.6155
$ws.Range("J71").Value = 3073.25
$ws.Range("L71").Value = 15366.25
$ws.Range("N71").Value = -22854.25
$ws.Range("H82").Value = 400.5
$ws.Range("I82").Value = 411.2
$ws.Range("K82").Value = 411.2
$ws.Range("M82").Value = -50.19999999999999
$ws.Range("H85").Value = 400.5
$ws.Range("I85").Value = 411.2
$ws.Range("K85").Value = 411.2
$ws.Range("M85").Value = 836.8
$ws.Range("H93").Value = 4287.1113
$ws.Range("J93").Value = 2902
$ws.Range("L93").Value = 2902
$ws.Range("N93").Value = -5398
$ws.Range("H126").Value = 3648.6
$ws.Range("I126").Value = 2248
$ws.Range("J126").Value = 3998.75
$ws.Range("K126").Value = 6744
$ws.Range("L126").Value = 11996.25
$ws.Range("M126").Value = -4274
$ws.Range("N126").Value = -16936.25
$ws.Range("H136").Value = 2238.48
$ws.Range("I136").Value = 2216.5
$ws.Range("J136").Value = 2399.6667
$ws.Range("K136").Value = 6649.5
$ws.Range("L136").Value = 7199.000100000001
$ws.Range("M136").Value = -4099.5
$ws.Range("N136").Value = -12299.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3509
$ws.Range("I4").Value = 3082
$ws.Range("J4").Value = 4149.5
$ws.Range("K4").Value = 3082
$ws.Range("L4").Value = 4149.5
$ws.Range("M4").Value = -2969
$ws.Range("N4").Value = -4375.5
